# correct $ -> \ in png #268
#
# Slide 4 has three places that describe the "less than 10,000 yen"
# code-list entry:
#   1) a Java "util:map" bean-definition sample (shape 1)
#   2) the rendered table that actually shows the mapping (shape 5)
#   3) an HTML <select>/<option> sample (shape 13)
#
# Samples (1) and (3) already read "...\10,000..." correctly but were
# split across several runs that all share identical formatting -- this
# change tidies those back into single runs. Sample (2) had a genuine
# typo ("$10,000" instead of "\10,000") which is the actual visible bug
# this commit fixes.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)

# ---------------------------------------------------------------------
# 1) Shape 1: Java bean-definition code sample, paragraph 7:
#      <entry key="10000" value="Less than \10,000" />
#    Keep the "<entry key=" run as-is (different formatting) and merge
#    the remaining three runs -- which already read "\10,000" -- into
#    one italic run.
# ---------------------------------------------------------------------
$shp1 = $s.Shapes.Item(1)
$tf1 = $shp1.TextFrame
$para1 = $tf1.TextRange.Paragraphs(7, 1)
$run1 = $para1.Runs(2, 1)
$start1 = $run1.Start
$len1 = $run1.Text.Length

# Re-assigning a TextRange's .Text to the exact text it already contains
# is treated as a no-op by the engine (the underlying runs are left
# untouched), so first overwrite the region with a throwaway placeholder
# of a different length. That forces the engine to rebuild the region as
# a single, freshly-formatted run.
$placeholder1 = "PLACEHOLDER_TEXT_FOR_RUN_MERGE_ONE"
$region1 = $tf1.TextRange.Characters($start1, $len1)
$region1.Text = $placeholder1

# Re-fetch the (now single, placeholder-filled) run and overwrite it
# with the real final text. Using the freshly read length keeps this
# correct even though PowerPoint's TextRange.Text includes a trailing
# carriage return when the range ends at a paragraph mark.
$placeholderRegion1 = $tf1.TextRange.Characters($start1, $placeholder1.Length)
$placeholderRegion1.Text = "`"10000`" value=`"Less than \10,000`" />"

# ---------------------------------------------------------------------
# 2) Shape 5: the rendered table. Row 3 / column 2 visibly (and
#    incorrectly) reads "Less than $10,000"; fix the stray "$" in place.
# ---------------------------------------------------------------------
$shp2 = $s.Shapes.Item(5)
$cell2 = $shp2.Table.Cell(3, 2)
$tf2 = $cell2.Shape.TextFrame
$cellText2 = $tf2.TextRange.Text
$dollarPos2 = $cellText2.IndexOf('$') + 1   # TextRange indices are 1-based
$dollarChar2 = $tf2.TextRange.Characters($dollarPos2, 1)
$dollarChar2.Text = "\"

# ---------------------------------------------------------------------
# 3) Shape 13: HTML <select>/<option> code sample, paragraph 4:
#      <option value="10000">Less than \10,000</option>
#    Merge the three runs (already reading "\10,000") into one run.
# ---------------------------------------------------------------------
$shp3 = $s.Shapes.Item(13)
$tf3 = $shp3.TextFrame
$para3 = $tf3.TextRange.Paragraphs(4, 1)
$start3 = $para3.Start
$len3 = $para3.Text.Length

$placeholder3 = "PLACEHOLDER_TEXT_FOR_RUN_MERGE_THREE_XYZ"
$region3 = $tf3.TextRange.Characters($start3, $len3)
$region3.Text = $placeholder3

$placeholderRegion3 = $tf3.TextRange.Characters($start3, $placeholder3.Length)
$placeholderRegion3.Text = "        <option value=`"10000`">Less than \10,000</option>"
